$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for new columns I and J, using same formatting as existing headers
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for columns I and J (identical values per row)
$values = @{
    2  = 7
    3  = 6
    4  = 6
    5  = 6
    6  = 8
    7  = 6
    8  = 7
    9  = 9
    10 = 7
    11 = 7
    12 = 9
    13 = 6
    14 = 7
    15 = 6
    16 = 8
    17 = 5
    18 = 8
    19 = 4
    20 = 5
}

foreach ($row in $values.Keys) {
    $v = $values[$row]
    $ws.Cells.Item($row, 9).Value = $v
    $ws.Cells.Item($row, 10).Value = $v
}
